$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Delete row 12 (未齊件說明 / NotYetItem) from the DBD table.
# This shifts rows 13-19 up to rows 12-18 and keeps their original
# content / shared formulas intact.
$ws.Rows.Item(12).Delete()

# The engine doesn't re-point the shared "=A{n-1}+1" formulas across the
# delete the way Excel does, so restore them explicitly.
$ws.Range("A12").Formula = "=A11+1"
$ws.Range("A13").Formula = "=A12+1"

# Restore the view that was saved with the workbook after the edit.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("I13").Select()
